$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column I ("Tipo Doc. Emisor"), shifting columns J:U left to I:T
$ws.Range("I1").EntireColumn.Delete()
